$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-10-30 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-31 Tuesday", 2) | Out-Null

# Update each multiplication problem in the table (old text values are unique
# within the document, so a simple find/replace targets the correct cell).
$d.Content.Find.Execute("66×66=", $true, $false, $false, $false, $false, $true, 1, $false, "95×93=", 2) | Out-Null
$d.Content.Find.Execute("12×50=", $true, $false, $false, $false, $false, $true, 1, $false, "46×66=", 2) | Out-Null
$d.Content.Find.Execute("53×17=", $true, $false, $false, $false, $false, $true, 1, $false, "24×68=", 2) | Out-Null
$d.Content.Find.Execute("16×30=", $true, $false, $false, $false, $false, $true, 1, $false, "28×46=", 2) | Out-Null
$d.Content.Find.Execute("22×36=", $true, $false, $false, $false, $false, $true, 1, $false, "53×78=", 2) | Out-Null
$d.Content.Find.Execute("94×47=", $true, $false, $false, $false, $false, $true, 1, $false, "79×83=", 2) | Out-Null
$d.Content.Find.Execute("39×38=", $true, $false, $false, $false, $false, $true, 1, $false, "80×39=", 2) | Out-Null
$d.Content.Find.Execute("27×44=", $true, $false, $false, $false, $false, $true, 1, $false, "90×72=", 2) | Out-Null
$d.Content.Find.Execute("38×86=", $true, $false, $false, $false, $false, $true, 1, $false, "85×56=", 2) | Out-Null
$d.Content.Find.Execute("11×58=", $true, $false, $false, $false, $false, $true, 1, $false, "38×36=", 2) | Out-Null
$d.Content.Find.Execute("28×31=", $true, $false, $false, $false, $false, $true, 1, $false, "30×75=", 2) | Out-Null
$d.Content.Find.Execute("31×81=", $true, $false, $false, $false, $false, $true, 1, $false, "92×21=", 2) | Out-Null
$d.Content.Find.Execute("46×58=", $true, $false, $false, $false, $false, $true, 1, $false, "79×17=", 2) | Out-Null
$d.Content.Find.Execute("99×61=", $true, $false, $false, $false, $false, $true, 1, $false, "64×23=", 2) | Out-Null
$d.Content.Find.Execute("92×82=", $true, $false, $false, $false, $false, $true, 1, $false, "49×53=", 2) | Out-Null
$d.Content.Find.Execute("22×41=", $true, $false, $false, $false, $false, $true, 1, $false, "87×42=", 2) | Out-Null
$d.Content.Find.Execute("74×22=", $true, $false, $false, $false, $false, $true, 1, $false, "89×40=", 2) | Out-Null
$d.Content.Find.Execute("66×59=", $true, $false, $false, $false, $false, $true, 1, $false, "51×85=", 2) | Out-Null
$d.Content.Find.Execute("17×64=", $true, $false, $false, $false, $false, $true, 1, $false, "14×31=", 2) | Out-Null
$d.Content.Find.Execute("85×54=", $true, $false, $false, $false, $false, $true, 1, $false, "55×78=", 2) | Out-Null
$d.Content.Find.Execute("23×46=", $true, $false, $false, $false, $false, $true, 1, $false, "48×28=", 2) | Out-Null
$d.Content.Find.Execute("44×59=", $true, $false, $false, $false, $false, $true, 1, $false, "19×97=", 2) | Out-Null
$d.Content.Find.Execute("51×34=", $true, $false, $false, $false, $false, $true, 1, $false, "35×98=", 2) | Out-Null
$d.Content.Find.Execute("46×83=", $true, $false, $false, $false, $false, $true, 1, $false, "89×13=", 2) | Out-Null
$d.Content.Find.Execute("75×43=", $true, $false, $false, $false, $false, $true, 1, $false, "35×54=", 2) | Out-Null
